$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Extracted on" timestamp in the subtitle cell.
$ws.Range("A2").Value = "This is an interesting study.Extracted on : 2022/09/26 12:12:12"

# Add an option for column renaming: rename the "label" header to "Variables"
# and clear the "levels" header (now blank, but still text-typed) in the
# summary table header row.
$ws.Range("A4").Value = "Variables"
$ws.Range("B4").Value = "'"
# Re-apply the original header formatting to B4 (the quote-prefix trick used
# above to keep the cell text-typed but empty also tags the style with a
# quote-prefix flag; restore the plain header style by copying it from a
# neighboring header cell that still carries it).
$ws.Range("C4").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false
